$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $newValue
    $c.Style = "Normal"
}

Set-TextValue "D2" "42.664.92"
Set-TextValue "E2" "  +1.68%  "
Set-TextValue "D3" "2.294.07"
Set-TextValue "E3" "  -0.10%  "
Set-TextValue "E4" "  +0.21%  "
Set-TextValue "D5" "322.48"
Set-TextValue "E5" "  +1.47%  "
Set-TextValue "D6" "104.35"
Set-TextValue "E6" "  +0.06%  "
Set-TextValue "E7" "  +0.31%  "
Set-TextValue "D9" "0.610"
Set-TextValue "E9" "  +0.51%  "
Set-TextValue "D10" "40.15"
Set-TextValue "E10" "  +2.24%  "
Set-TextValue "E11" "  +0.00%  "
Set-TextValue "D12" "8.54"
Set-TextValue "E12" "  +2.14%  "
Set-TextValue "E13" "  -0.10%  "
Set-TextValue "D14" "0.971"
Set-TextValue "D15" "15.23"
Set-TextValue "E15" "  -0.49%  "
Set-TextValue "D16" "2.642.74"
Set-TextValue "E16" "  -0.05%  "
Set-TextValue "D17" "2.286.66"
Set-TextValue "E17" "  +0.28%  "
Set-TextValue "D18" "42.600.54"
Set-TextValue "E18" "  +1.54%  "
Set-TextValue "D19" "7.45"
Set-TextValue "E19" "  -2.58%  "
Set-TextValue "E20" "  +0.08%  "
Set-TextValue "D21" "13.34"
Set-TextValue "E21" "  +32.84%  "
Set-TextValue "D22" "73.27"
Set-TextValue "E22" "  -0.59%  "
Set-TextValue "E23" "  +0.16%  "
Set-TextValue "D24" "270.14"
Set-TextValue "E24" "  -6.64%  "
Set-TextValue "E25" "  -1.90%  "
Set-TextValue "E26" "  -0.37%  "
Set-TextValue "D27" "10.90"
Set-TextValue "E27" "  +0.18%  "
Set-TextValue "E28" "  +2.91%  "
Set-TextValue "D29" "38.38"
Set-TextValue "E29" "  +9.19%  "
Set-TextValue "D30" "22.56"
Set-TextValue "E30" "  -3.12%  "
Set-TextValue "D31" "165.55"
Set-TextValue "E31" "  +1.11%  "
Set-TextValue "D32" "6.16"
Set-TextValue "E32" "  +5.15%  "
Set-TextValue "D33" "0.0880"
Set-TextValue "E33" "  -0.23%  "
Set-TextValue "E34" "  +0.38%  "
Set-TextValue "D35" "0.114"
Set-TextValue "E35" "  -1.43%  "
Set-TextValue "E36" "  -13.87%  "
Set-TextValue "D37" "4.64"
Set-TextValue "D38" "0.0356"
Set-TextValue "E38" "  +1.71%  "
Set-TextValue "D39" "3.74"
Set-TextValue "E39" "  +3.54%  "
Set-TextValue "D40" "2.73"
Set-TextValue "E40" "  -4.24%  "
Set-TextValue "E41" "  +4.11%  "
Set-TextValue "D42" "70.06"
Set-TextValue "E42" "  -1.11%  "
Set-TextValue "E43" "  +0.22%  "
Set-TextValue "E44" "  -0.17%  "
Set-TextValue "D45" "93.37"
Set-TextValue "E45" "  -9.46%  "
Set-TextValue "D46" "12.31"
Set-TextValue "E46" "  +2.11%  "
Set-TextValue "D47" "81.49"
Set-TextValue "E47" "  +4.09%  "
Set-TextValue "D48" "113.65"
Set-TextValue "E48" "  -1.99%  "
Set-TextValue "D49" "8.91"
Set-TextValue "E49" "  -1.91%  "
Set-TextValue "D50" "5.27"
Set-TextValue "E50" "  -1.13%  "
Set-TextValue "D51" "1.585.84"
Set-TextValue "E51" "  +1.85%  "
